# Fabrikam Q1 marketing campaigns - header row updates
# - Rename "Fecha de inicio" -> "Fecha de lanzamiento"
# - Rename "Usuarios objetivo totales" -> "Total de usuarios seleccionados"
# - Bold every header cell in row 1 (A1:H1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# White color used by the header font (rgb FFFFFFFF -> RGB(255,255,255))
$White = 16777215

function Set-HeaderBold($cell) {
    $len = $cell.Characters().Count()
    if ($len -gt 1) {
        $cell.Characters(1, $len - 1).Font.Bold = $true
        $cell.Characters($len, 1).Font.Bold = $true
    } else {
        $cell.Characters(1, $len).Font.Bold = $true
    }
}

function Set-HeaderWhite($cell) {
    $len = $cell.Characters().Count()
    if ($len -gt 1) {
        $cell.Characters(1, $len - 1).Font.Color = $White
        $cell.Characters($len, 1).Font.Color = $White
    } else {
        $cell.Characters(1, $len).Font.Color = $White
    }
}

# Rename the two headers first (new text inherits the default/black font,
# so we re-apply the white color explicitly afterwards).
$ws.Range("C1").Value = "Fecha de lanzamiento"
$ws.Range("G1").Value = "Total de usuarios seleccionados"

Set-HeaderWhite($ws.Range("C1"))
Set-HeaderWhite($ws.Range("G1"))

# Bold every header cell (A1, C1, D1, G1, H1 were not bold before;
# B1, E1, F1 were already bold).
Set-HeaderBold($ws.Range("A1"))
Set-HeaderBold($ws.Range("C1"))
Set-HeaderBold($ws.Range("D1"))
Set-HeaderBold($ws.Range("G1"))
Set-HeaderBold($ws.Range("H1"))
